$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "28.874.71"
$ws.Range("E2").Value2 = "  +8.00%  "
$ws.Range("D3").Value2 = "1.814.84"
$ws.Range("E3").Value2 = "  +5.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value2 = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.82"
$ws.Range("E5").Value2 = "  +2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value2 = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4922"
$ws.Range("E7").Value2 = "  +2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.74"
$ws.Range("E8").Value2 = "  +6.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2771"
$ws.Range("E9").Value2 = "  +7.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06387"
$ws.Range("D11").Value2 = "1.808.09"
$ws.Range("E11").Value2 = "  +4.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.72"
$ws.Range("E12").Value2 = "  +5.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07073"
$ws.Range("E13").Value2 = "  +3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6418"
$ws.Range("E14").Value2 = "  +6.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.88"
$ws.Range("E15").Value2 = "  +9.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.660"
$ws.Range("E16").Value2 = "  +4.58%  "
$ws.Range("D17").Value2 = "28.896.97"
$ws.Range("E17").Value2 = "  +8.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value2 = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007307"
$ws.Range("E19").Value2 = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value2 = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value2 = "  +7.50%  "
$ws.Range("D22").Value2 = "2.040.11"
$ws.Range("E22").Value2 = "  +4.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.548"
$ws.Range("E23").Value2 = "  +3.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.763"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.350"
$ws.Range("E25").Value2 = "  +5.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.64"
$ws.Range("E26").Value2 = "  +3.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "128.52"
$ws.Range("E27").Value2 = "  +21.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.35"
$ws.Range("E28").Value2 = "  +7.58%  "
$ws.Range("E29").Value2 = "  +6.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value2 = "  +2.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.121"
$ws.Range("E31").Value2 = "  +2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08334"
$ws.Range("E32").Value2 = "  +5.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.770"
$ws.Range("E33").Value2 = "  +3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04906"
$ws.Range("E34").Value2 = "  +8.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.094"
$ws.Range("E35").Value2 = "  +9.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.696"
$ws.Range("E36").Value2 = "  +3.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6700"
$ws.Range("E37").Value2 = "  +8.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.282"
$ws.Range("E38").Value2 = "  +14.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.677"
$ws.Range("E39").Value2 = "  +9.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9505"
$ws.Range("E40").Value2 = "  +2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.147"
$ws.Range("E41").Value2 = "  +9.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01583"
$ws.Range("E42").Value2 = "  +5.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value2 = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.49"
$ws.Range("E44").Value2 = "  +0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4065"
$ws.Range("E45").Value2 = "  +6.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.129"
$ws.Range("E46").Value2 = "  +5.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1219"
$ws.Range("E47").Value2 = "  +5.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05512"
$ws.Range("E48").Value2 = "  +2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.080"
$ws.Range("E49").Value2 = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.57"
$ws.Range("E50").Value2 = "  +5.01%  "
$ws.Range("B51").Value2 = "NEARProtocol"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.303"
$ws.Range("E51").Value2 = "  +4.48%  "
